# Scheduled data refresh: update cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the ALC, ARM, BSM,
# CRP, CUL, GSM and LTW sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4233.476
$ws.Range("J17").Value = 4233.476
$ws.Range("L17").Value = 12700.428
$ws.Range("N17").Value = -13036.428

$ws.Range("H41").Value = 316.4
$ws.Range("I41").Value = 394.73334
$ws.Range("J41").Value = 198.9
$ws.Range("K41").Value = 394.73334
$ws.Range("L41").Value = 198.9
$ws.Range("M41").Value = 45.26666
$ws.Range("N41").Value = -1078.9

$ws.Range("H53").Value = 470.85715
$ws.Range("I53").Value = 280.7647
$ws.Range("K53").Value = 280.7647
$ws.Range("M53").Value = 356.2353

$ws.Range("H129").Value = 1945.1875
$ws.Range("I129").Value = 1466.4546
$ws.Range("J129").Value = 2998.4
$ws.Range("K129").Value = 4399.3638
$ws.Range("L129").Value = 8995.200000000001
$ws.Range("M129").Value = 600.6361999999999
$ws.Range("N129").Value = -18995.2

$ws.Range("H132").Value = 1067.7307
$ws.Range("I132").Value = 858.3261
$ws.Range("K132").Value = 2574.9783
$ws.Range("M132").Value = -44.97829999999976

$ws.Range("H137").Value = 4749.073
$ws.Range("I137").Value = 2034.421
$ws.Range("J137").Value = 7093.5454
$ws.Range("K137").Value = 6103.263
$ws.Range("L137").Value = 21280.6362
$ws.Range("M137").Value = -3553.263
$ws.Range("N137").Value = -26380.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19605.012
$ws.Range("I32").Value = 19552.041
$ws.Range("J32").Value = 20249.5
$ws.Range("K32").Value = 19552.041
$ws.Range("L32").Value = 20249.5
$ws.Range("M32").Value = -19265.041
$ws.Range("N32").Value = -20823.5

$ws.Range("H61").Value = 7227.032
$ws.Range("I61").Value = 5062.5264
$ws.Range("J61").Value = 10654.167
$ws.Range("K61").Value = 5062.5264
$ws.Range("L61").Value = 10654.167
$ws.Range("M61").Value = -4850.5264
$ws.Range("N61").Value = -11078.167

$ws.Range("H74").Value = 280911.22
$ws.Range("I74").Value = 304327.44
$ws.Range("J74").Value = 23333
$ws.Range("K74").Value = 304327.44
$ws.Range("L74").Value = 23333
$ws.Range("M74").Value = -303453.44
$ws.Range("N74").Value = -25081

$ws.Range("H77").Value = 280911.22
$ws.Range("I77").Value = 304327.44
$ws.Range("J77").Value = 23333
$ws.Range("K77").Value = 1521637.2
$ws.Range("L77").Value = 116665
$ws.Range("M77").Value = -1517269.2
$ws.Range("N77").Value = -125401

$ws.Range("H132").Value = 7212.447
$ws.Range("I132").Value = 5366.811
$ws.Range("K132").Value = 16100.433
$ws.Range("M132").Value = -13570.433

$ws.Range("H136").Value = 7227.032
$ws.Range("I136").Value = 5062.5264
$ws.Range("J136").Value = 10654.167
$ws.Range("K136").Value = 15187.5792
$ws.Range("L136").Value = 31962.501
$ws.Range("M136").Value = -12637.5792
$ws.Range("N136").Value = -37062.501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1477.5
$ws.Range("I107").Value = 1231
$ws.Range("K107").Value = 1231
$ws.Range("M107").Value = 689

$ws.Range("H134").Value = 4532.323
$ws.Range("I134").Value = 3173.9138
$ws.Range("K134").Value = 9521.741399999999
$ws.Range("M134").Value = -6986.741399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23261820
$ws.Range("I31").Value = 58826844
$ws.Range("J31").Value = 7767.8076
$ws.Range("K31").Value = 58826844
$ws.Range("L31").Value = 7767.8076
$ws.Range("M31").Value = -58826549
$ws.Range("N31").Value = -8357.8076

$ws.Range("H34").Value = 23261820
$ws.Range("I34").Value = 58826844
$ws.Range("J34").Value = 7767.8076
$ws.Range("K34").Value = 58826844
$ws.Range("L34").Value = 7767.8076
$ws.Range("M34").Value = -58826642
$ws.Range("N34").Value = -8171.8076

$ws.Range("H58").Value = 4231
$ws.Range("I58").Value = 1903.9131
$ws.Range("K58").Value = 1903.9131
$ws.Range("M58").Value = -1700.9131

$ws.Range("H62").Value = 17666.6
$ws.Range("I62").Value = 13998
$ws.Range("J62").Value = 20112.334
$ws.Range("K62").Value = 13998
$ws.Range("L62").Value = 20112.334
$ws.Range("M62").Value = -13374
$ws.Range("N62").Value = -21360.334

$ws.Range("H65").Value = 17666.6
$ws.Range("I65").Value = 13998
$ws.Range("J65").Value = 20112.334
$ws.Range("K65").Value = 69990
$ws.Range("L65").Value = 100561.67
$ws.Range("M65").Value = -66870
$ws.Range("N65").Value = -106801.67

$ws.Range("H107").Value = 967.875
$ws.Range("I107").Value = 967.875
$ws.Range("K107").Value = 967.875
$ws.Range("M107").Value = 952.125

$ws.Range("H136").Value = 4231
$ws.Range("I136").Value = 1903.9131
$ws.Range("K136").Value = 5711.7393
$ws.Range("M136").Value = -3161.7393

$ws.Range("H141").Value = 192290.7
$ws.Range("J141").Value = 221348.88
$ws.Range("L141").Value = 221348.88
$ws.Range("N141").Value = -231708.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws.Range("H129").Value = 27779590
$ws.Range("I129").Value = 643.5
$ws.Range("J129").Value = 62503270
$ws.Range("K129").Value = 1930.5
$ws.Range("L129").Value = 187509810
$ws.Range("M129").Value = 3069.5
$ws.Range("N129").Value = -187519810

$ws.Range("H131").Value = 11500150
$ws.Range("I131").Value = 23810694
$ws.Range("J131").Value = 10308.134
$ws.Range("K131").Value = 71432082
$ws.Range("L131").Value = 30924.402
$ws.Range("M131").Value = -71427042
$ws.Range("N131").Value = -41004.402

$ws.Range("H137").Value = 8963.857
$ws.Range("I137").Value = 21199.2
$ws.Range("J137").Value = 2166.4443
$ws.Range("K137").Value = 63597.60000000001
$ws.Range("L137").Value = 6499.3329
$ws.Range("M137").Value = -58497.60000000001
$ws.Range("N137").Value = -16699.3329

$ws.Range("H140").Value = 27779142
$ws.Range("I140").Value = 50001180
$ws.Range("J140").Value = 1593.75
$ws.Range("K140").Value = 150003540
$ws.Range("L140").Value = 4781.25
$ws.Range("M140").Value = -149998360
$ws.Range("N140").Value = -15141.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1991.7142
$ws.Range("J102").Value = 4332.6665
$ws.Range("L102").Value = 4332.6665
$ws.Range("N102").Value = -7576.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3418.3784
$ws.Range("I22").Value = 1926
$ws.Range("K22").Value = 1926
$ws.Range("M22").Value = -1631

$ws.Range("H27").Value = 3418.3784
$ws.Range("I27").Value = 1926
$ws.Range("K27").Value = 1926
$ws.Range("M27").Value = -1819

$ws.Range("H82").Value = 2787.375
$ws.Range("I82").Value = 2079.8
$ws.Range("J82").Value = 3966.6667
$ws.Range("K82").Value = 2079.8
$ws.Range("L82").Value = 3966.6667
$ws.Range("M82").Value = -1718.8
$ws.Range("N82").Value = -4688.6667

$ws.Range("H85").Value = 2787.375
$ws.Range("I85").Value = 2079.8
$ws.Range("J85").Value = 3966.6667
$ws.Range("K85").Value = 2079.8
$ws.Range("L85").Value = 3966.6667
$ws.Range("M85").Value = -831.8000000000002
$ws.Range("N85").Value = -6462.6667

$ws.Range("H132").Value = 5799.161
$ws.Range("I132").Value = 5129.4907
$ws.Range("J132").Value = 6950.1562
$ws.Range("K132").Value = 15388.4721
$ws.Range("L132").Value = 20850.4686
$ws.Range("M132").Value = -12858.4721
$ws.Range("N132").Value = -25910.4686
